$d = $word.ActiveDocument

# 1) "Sistema exibira..." popup description: desse -> deste, add Fornecedor
$d.Content.Find.Execute(
    "Sistema exibirá um popup com todas as informações desse produto: Código, Nome, Tipo e Valor.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Sistema exibirá um popup com todas as informações deste produto: Código, Nome, Tipo, Fornecedor e Valor.",
    2) | Out-Null

# 2) "Ator e advertido..." -> "Ator e notificado com um popup..."
$d.Content.Find.Execute(
    "2. Ator é advertido sobre a exclusão do produto;",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "2. Ator é notificado com um popup sobre a exclusão do produto;",
    2) | Out-Null

# 3) Historico table: split last row's "Alteracao" text into two lines,
#    then append a new history row for the current update.
$t = $d.Tables.Item($d.Tables.Count)
$lastRow = $t.Rows.Item($t.Rows.Count)
$descCell = $lastRow.Cells.Item(3)
$descCell.Range.Text = "Atualização da Descrição do " + [char]13 + "Caso de Uso"

$newRow = $t.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "07/12/2020"
$newRow.Cells.Item(2).Range.Text = "Wagner Prata"
$newRow.Cells.Item(3).Range.Text = "Atualização da Descrição do " + [char]13 + "Caso de Uso"

# 4) Insert two new empty paragraphs right after the table (before the
#    existing trailing empty paragraph). Use the document's content end
#    (minus the final paragraph mark) as the anchor so this keeps working
#    even though the table mutations above can leave stale Paragraphs
#    collections behind.
$endPos = $d.Content.End
$rng = $d.Range($endPos - 1, $endPos - 1)
$rng.InsertParagraphBefore()

$endPos2 = $d.Content.End
$rng2 = $d.Range($endPos2 - 1, $endPos2 - 1)
$rng2.InsertParagraphBefore()
